$d = $word.ActiveDocument

# The document contains two "2024.8.21   天气晴" date-heading paragraphs
# (an accidental duplicate day). We must only touch the second one -
# the one immediately followed by the "...第三次课程..." paragraph - and
# change its date to "2024.8.22", splitting the run so a "_GoBack"
# bookmark sits right after the date and before the "   天气晴" suffix.
# The stray "_GoBack" bookmark that currently sits at the very end of the
# "...第三次课程..." paragraph must be removed (it effectively relocates
# there, since bookmark names are unique).

$paras = $d.Paragraphs
$count = $paras.Count

$targetIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $ptext = $paras.Item($i).Range.Text
    if ($ptext -like "*2024.8.21*天气晴*") {
        if ($i -lt $count) {
            $nextText = $paras.Item($i + 1).Range.Text
            if ($nextText -like "*第三次课程*") {
                $targetIndex = $i
            }
        }
    }
}

$targetPara = $paras.Item($targetIndex)
$targetRange = $targetPara.Range

# Replace "2024.8.21" with "2024.8.22" within just this paragraph,
# keeping the rest of the run ("   天气晴") untouched so formatting stays
# on a single run for now.
$null = $targetRange.Find.Execute("2024.8.21", $true, $false, $false, $false, $false, `
                                   $true, 1, $false, "2024.8.22", 2)

# Re-fetch the paragraph range (Find may have adjusted bounds) and compute
# the split point right after "2024.8.22" (9 characters long), where the
# bookmark must be inserted, splitting the single run into two runs with
# identical formatting.
$targetPara = $paras.Item($targetIndex)
$splitPos = $targetPara.Range.Start + 9
$splitRange = $d.Range($splitPos, $splitPos)

# Adding a bookmark with a name that already exists simply relocates it,
# which both inserts it at the new spot and removes it from its old one.
$null = $d.Bookmarks.Add("_GoBack", $splitRange)
